# Updated multiple management files to have consistent and english values.
#
# On the "Bon_2017" sheet, the fertilizer amount columns (D) held the amount
# as text like "50 kg/ha" and the notice columns (E) referenced the German
# fertilizer name "N-Düngung KAS". These are normalized to a plain numeric
# amount (kg/ha implied by context) and the internationally recognised
# fertilizer abbreviation "CAN" (Calcium Ammonium Nitrate).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bon_2017")

$rows = @(55, 56, 57, 66, 67, 68)
$amounts = @{ 55 = 50; 56 = 25; 57 = 60; 66 = 50; 67 = 25; 68 = 60 }

foreach ($r in $rows) {
    $ws.Range("D$r").Value = $amounts[$r]
    $ws.Range("E$r").Value = "CAN"
}
